# Fixed update to excel issue
# 1. Rename header "Requested quantity" -> "Weekly_PO_Qty" on "Weekly Quantity" sheet
# 2. Rename header "Requested quantity" -> "Monthly_PO_Qty" on "Monthly Trend" sheet
# 3. Add a new "PO Forecast" sheet with forecast data (ds, PO_Forecast, yhat_lower, yhat_upper)

$wb = $excel.ActiveWorkbook

# --- 1. Weekly Quantity sheet: rename header cell B1 ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# --- 2. Monthly Trend sheet: rename header cell B1 ---
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 3. Add new "PO Forecast" sheet at the end ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Match the header formatting used on the other two sheets (bold, centered, thin border)
$headerRange = $wsForecast.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Forecast data rows: ds, PO_Forecast, yhat_lower, yhat_upper
$data = @(
    @(45060.99999999999,75,-58.19987665253871,212.4864789086359),
    @(45067.99999999999,75,-64.9494588861117,215.2851665610637),
    @(45074.99999999999,75,-61.98435425811846,210.8630939228354),
    @(45081.99999999999,75,-62.33552229011428,209.163796893517),
    @(45088.99999999999,75,-63.44628683802269,208.6358455470238),
    @(45095.99999999999,75,-64.35842268641335,218.6634555934569),
    @(45102.99999999999,75,-67.35706522849576,213.2142790590892),
    @(45109.99999999999,75,-76.03931392406344,206.6729678873885),
    @(45116.99999999999,75,-69.66522955206125,213.5634066529048),
    @(45123.99999999999,75,-69.11732962598512,209.115069965591),
    @(45130.99999999999,75,-57.74299915966228,222.360474232781),
    @(45137.99999999999,75,-69.8281811329943,211.619808120141),
    @(45144.99999999999,75,-62.90624656853009,211.9862119407166),
    @(45151.99999999999,75,-68.71275014065922,219.0923507080149),
    @(45158.99999999999,75,-66.42918282883309,220.6092580457613),
    @(45193.99999999999,75,-77.44716237730496,205.0328481612531),
    @(45214.99999999999,75,-73.29696722842719,205.7385523622323),
    @(45221.99999999999,75,-66.46635055999741,215.4154799582075),
    @(45228.99999999999,75,-70.00205756331449,225.4388018237953),
    @(45235.99999999999,75,-54.7877683945605,219.1398591116272),
    @(45242.99999999999,75,-76.01112591628441,205.6521046762794),
    @(45249.99999999999,75,-53.50912831771464,210.8024236583533),
    @(45256.99999999999,75,-59.75188065444852,223.0986859373758),
    @(45263.99999999999,75,-65.90460653184118,216.2901063179149),
    @(45270.99999999999,75,-65.89600589739975,208.9880055962775),
    @(45277.99999999999,75,-63.82103878587405,217.3592062114963),
    @(45354.99999999999,75,-74.27788081587551,214.5610162845635),
    @(45389.99999999999,75,-64.23320682892812,210.5387413766938),
    @(45403.99999999999,75,-71.5672659651824,220.0231968377726),
    @(45410.99999999999,75,-64.90583331360511,207.1740056375054),
    @(45424.99999999999,75,-68.81870031372529,219.3910470959515),
    @(45431.99999999999,75,-66.51655367563471,215.708776424036),
    @(45438.99999999999,75,-74.720585499547,209.3292727432034),
    @(45452.99999999999,75,-56.04594229582679,220.6752314456259),
    @(45459.99999999999,75,-61.20202652334665,215.3090301784703),
    @(45466.99999999999,75,-63.07707782553893,222.8532721446041),
    @(45473.99999999999,75,-77.72811190516843,206.7158990211975),
    @(45480.99999999999,75,-76.9216938177978,212.2532264535268),
    @(45487.99999999999,75,-70.31965270716333,213.0833912937699),
    @(45494.99999999999,75,-67.50642323043424,215.1837646869265),
    @(45550.99999999999,74,-68.14314022055875,217.4509612827292),
    @(45557.99999999999,74,-52.82965900091621,210.8270335520784),
    @(45564.99999999999,74,-65.99097184254235,206.9030574651449),
    @(45571.99999999999,74,-64.46195863629062,218.8777939908945),
    @(45578.99999999999,74,-69.71956322740817,219.411803902972),
    @(45585.99999999999,75,-60.37801682175088,215.5123689673633),
    @(45592.99999999999,74,-68.3279847729157,203.6712669065625),
    @(45599.99999999999,74,-68.7319774626271,204.851288995486),
    @(45606.99999999999,74,-64.89619132276059,219.0764935261596)
)

$row = 2
foreach ($rec in $data) {
    $wsForecast.Cells.Item($row, 1).Value = $rec[0]
    $wsForecast.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $wsForecast.Cells.Item($row, 2).Value = $rec[1]
    $wsForecast.Cells.Item($row, 3).Value = $rec[2]
    $wsForecast.Cells.Item($row, 4).Value = $rec[3]
    $row++
}
